$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F7: update text only (style s=15 already set)
$ws.Range("F7").Value2 = "ripresa bayes rule (marginal, joint; fatto risolvere esempio infected | positive;  goat)"

# F8: update text, and apply the same formatting as D10 (style s=3)
# (new shared string #39, matches the index ordering of the target workbook)
$ws.Range("D10").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F8").Value2 = "beta-binomial (beta; con esercizio di tune prior  e inizio 3.9)"

# F16: new cell, new style -> yellow fill, centered text
# (new shared string #40)
$ws.Range("F16").Value2 = "presentare assignment e progetto"
$ws.Range("F16").Interior.Color = 65535
$ws.Range("F16").HorizontalAlignment = -4108

# E10: clear the stray content (moved to F10)
$ws.Range("E10").ClearContents()

# F10: new cell, formatted like D10 (style s=3)
# (new shared string #41)
$ws.Range("D10").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Range("F10").Value2 = "ripresa beta; binomial likelihood; up to slide 58 (sensitivity to the prior)"

# F14: new cell, formatted like B3 (style s=16, italic font)
# (new shared string #42)
$ws.Range("B3").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").Value2 = "gc: guardare come inserire test ipotesi di francesca"

# F9: new cell, formatted like D9 (style s=11), reuses existing shared string #7
$ws.Range("D9").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F9").Value2 = "exe MF: bayes rule"

# F11: new cell, formatted like D10 (style s=3), reuses existing shared string #9
$ws.Range("D10").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F11").Value2 = "finire beta-bin; 3.9, 3.10, normal-normal"

# F12: new cell, formatted like D10 (style s=3), reuses existing shared string #10
$ws.Range("D10").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$ws.Range("F12").Value2 = "reporting analysis; normal normal (esercizio likelihood)"

$excel.CutCopyMode = 0

# Update selection to F14 as in the final workbook
$ws.Range("F14").Select()
